$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (index 11)
$ws.Range("B7").Value = 0.2603692549741538
$ws.Range("C7").Value = 1.880313879190301
$ws.Range("D7").Value = 12.64824988312696
$ws.Range("E7").Value = 3.556437808134279
$ws.Range("F7").Value = 3.597207164417901
$ws.Range("G7").Value = 36

# Row 8 (index 12)
$ws.Range("B8").Value = 0.2149774803876897
$ws.Range("C8").Value = 2.066399049159351
$ws.Range("D8").Value = 12.86774599829277
$ws.Range("E8").Value = 3.58716406068816
$ws.Range("F8").Value = 3.632992502715109
$ws.Range("G8").Value = 35

# Row 9 (index 13)
$ws.Range("B9").Value = 0.3416184366656868
$ws.Range("C9").Value = 2.567298676063721
$ws.Range("D9").Value = 20.59292190581267
$ws.Range("E9").Value = 4.53794247493428
$ws.Range("F9").Value = 4.642619354761157
$ws.Range("G9").Value = 20

# Row 10 (index 14)
$ws.Range("B10").Value = -0.7803977842537473
$ws.Range("C10").Value = 2.218913924505684
$ws.Range("D10").Value = 9.291300011801141
$ws.Range("E10").Value = 3.048163383383696
$ws.Range("F10").Value = 3.066888094792081
$ws.Range("G10").Value = 13

# Row 11 (index 15)
$ws.Range("B11").Value = -0.6077611373159695
$ws.Range("C11").Value = 2.229470670271877
$ws.Range("D11").Value = 6.910496118181399
$ws.Range("E11").Value = 2.628782250050658
$ws.Range("F11").Value = 2.85944105511676
$ws.Range("G11").Value = 5
